# Insert a new data row before the current row 120 (La Araucanía / Berenjena
# weekly price log), pushing all subsequent rows down by one, and fill the
# new row with the latest week's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("120:120").Insert()

$ws.Cells.Item(120, 1).Value  = 10
$ws.Cells.Item(120, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(120, 3).Value  = "La Araucanía"
$ws.Cells.Item(120, 4).Value  = 44589
$ws.Cells.Item(120, 5).Value  = 9
$ws.Cells.Item(120, 6).Value  = 100112001
$ws.Cells.Item(120, 7).Value  = "Berenjena"
$ws.Cells.Item(120, 8).Value  = "Sin especificar"
$ws.Cells.Item(120, 9).Value  = "Primera"
$ws.Cells.Item(120, 10).Value = 55
$ws.Cells.Item(120, 11).Value = 10000
$ws.Cells.Item(120, 12).Value = 12000
$ws.Cells.Item(120, 13).Value = 11091
$ws.Cells.Item(120, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(120, 15).Value = "Región del Maule"
$ws.Cells.Item(120, 16).Value = 185
$ws.Cells.Item(120, 17).Value = 60
$ws.Cells.Item(120, 18).Value = "Hortaliza"
